$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The two draft rows (101/102) describing "nudgincsm" and "nudgincswe" move
# up into the main table (rows 67/68) with corrected / expanded content.
# Remove the old placeholder rows first.
# ---------------------------------------------------------------------------
$ws.Rows.Item(101).Delete()
$ws.Rows.Item(101).Delete()

# ---------------------------------------------------------------------------
# Row 67: nudgincsm
# ---------------------------------------------------------------------------
$ws.Range("A67").Value = "Eday"
$ws.Range("B67").Value = "nudgincsm"

$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "1"
$ws.Range("C67").NumberFormat = "General"

$ws.Range("D67").Value = "longitude latitude time"
$ws.Range("E67").Value = "Nudging Increment of Water in Soil Moisture"
$ws.Range("F67").Value = "kg m-2"
$ws.Range("G67").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/01c8c41a-a0d8-11e6-bc63-ac72891c3257.html","web")'

$ws.Range("H67").Value = "To be implemented:  grib 126.151:  ifs code name = 151.126  part of MFPPHY.  Have to be  made available via PEXTRA, upto now with some  non-defined or adhoc grib code. Nudincsm is, consistent with sm, saved for each of the four soil layers"
$ws.Range("H67").Characters(1, 19).Font.Color = 0
$ws.Range("H67").Characters(20, 13).Font.Color = 1972430
$ws.Range("H67").Characters(33, 207).Font.Color = 0

$ws.Range("I67").Value = "Emanuel Dutra, Wilhelm May, Thomas Reerink"
$ws.Range("J67").Value = "A nudging increment refers to an amount added to parts of a model system. The phrase 'nudging_increment_in_X' refers to an increment in quantity X over a time period which should be defined in the bounds of the time coordinate. 'Content' indicates a quantity per unit area. 'Water' means water in all phases. The mass content of water in soil refers to the vertical integral from the surface down to the bottom of the soil model. The 'soil content' of a quantity refers to the vertical integral from the surface down to the bottom of the soil model. For the content between specified levels in the soil, standard names including 'content_of_soil_layer' are used."
$ws.Range("K67").Value = "LS3MIP"

# ---------------------------------------------------------------------------
# Row 68: nudgincswe
# ---------------------------------------------------------------------------
$ws.Range("A68").Value = "Eday"
$ws.Range("B68").Value = "nudgincswe"

$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "1"
$ws.Range("C68").NumberFormat = "General"

$ws.Range("D68").Value = "longitude latitude time"
$ws.Range("E68").Value = "Nudging Increment of Water in Snow"
$ws.Range("F68").Value = "kg m-2"
$ws.Range("G68").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/0abbdddc-a0d8-11e6-bc63-ac72891c3257.html","web")'

$ws.Range("H68").Value = "To be implemented:  grib 126.152:  ifs code name = 152.126  part of MFPPHY.  Have to be  made available via PEXTRA, upto now with some  non-defined or adhoc grib code."
$ws.Range("H68").Characters(1, 20).Font.Color = 0
$ws.Range("H68").Characters(21, 12).Font.Color = 1972430
$ws.Range("H68").Characters(33, 135).Font.Color = 0

$ws.Range("I68").Value = "Emanuel Dutra, Wilhelm May, Thomas Reerink"
$ws.Range("J68").Value = "A nudging increment refers to an amount added to parts of a model system. The phrase 'nudging_increment_in_X' refers to an increment in quantity X over a time period which should be defined in the bounds of the time coordinate. The surface called 'surface' means the lower boundary of the atmosphere. 'Amount' means mass per unit area. 'Snow and ice on land' means ice in glaciers, ice caps, ice sheets & shelves, river and lake ice, any other ice on a land surface, such as frozen flood water, and snow lying on such ice or on the land surface."
$ws.Range("K68").Value = "LS3MIP"

# ---------------------------------------------------------------------------
# Update the view: selection / top-left cell moved near the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A67").Select()
